# Apply edits to AFFcon.xlsx (Sheet1):
#  - Swap each value in column F (rows 2-81) between "k" and "l"
#  - Change the sheet view: remove the A64 top-left scroll position and
#    move the selection to I95 (single cell) instead of F2:F81

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 81; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F
    $cur = $cell.Value2
    if ($cur -eq "k") {
        $cell.Value2 = "l"
    } elseif ($cur -eq "l") {
        $cell.Value2 = "k"
    }
}

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I95").Select()
